# Update Moldova MSME summary figures with more precise (two-decimal) values.
# These cells hold numeric-looking values that are stored as text in the
# workbook, so we force the number format to Text ("@") before assigning the
# new value - otherwise Excel would auto-convert the string into a real
# number, which is not what the source data intends.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Enterprises density (per 1000 people)  -> Micro / SMEs / MSMEs
$ws.Range("B11:D11").NumberFormat = "@"
$ws.Range("B11").Value = "10.77"
$ws.Range("C11").Value = "3.12"
$ws.Range("D11").Value = "13.89"

# Row 12: Employment (% of total) -> Micro / SMEs / MSMEs
$ws.Range("B12:D12").NumberFormat = "@"
$ws.Range("B12").Value = "17.39"
$ws.Range("C12").Value = "40.35"
$ws.Range("D12").Value = "57.75"

# Row 14: Enterprises (% of total) -> Micro / SMEs / MSMEs
$ws.Range("B14:D14").NumberFormat = "@"
$ws.Range("B14").Value = "75.64"
$ws.Range("C14").Value = "21.92"
$ws.Range("D14").Value = "97.56"
